# Replace the old "row 1 / row 2" gamelog data (which pulled text from
# xl/sharedStrings.xml) with the new, wider rows that add player identity
# columns (Last name, First name, Position) up front and a trailing
# fantasy-points numeric column, per the authoring script's rewrite of
# 2018_data/Keenan_Reynolds2018.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to stay text (so things like "2018-09-17" or "2" or
    # "24.278" are not silently reinterpreted by Excel as a date/number),
    # then drop the temporary "@" number-format again so we don't leave
    # stray formatting behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---- Row 1: Keenan Reynolds, WR, week 2 @ CHI (loss) ----
$ws.Range("A1").Value = "Reynolds"
$ws.Range("B1").Value = "Keenan"
$ws.Range("C1").Value = "WR"
Set-TextValue $ws.Range("D1") "2018-09-17"
Set-TextValue $ws.Range("E1") "2"
Set-TextValue $ws.Range("F1") "24.278"
$ws.Range("G1").Value = "SEA"
$ws.Range("H1").Value = "@"
$ws.Range("I1").Value = "CHI"
$ws.Range("J1").Value = "L 17-24"
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = 0

# ---- Row 2: Keenan Reynolds, WR, week 3 vs DAL (win) ----
$ws.Range("A2").Value = "Reynolds"
$ws.Range("B2").Value = "Keenan"
$ws.Range("C2").Value = "WR"
Set-TextValue $ws.Range("D2") "2018-09-23"
Set-TextValue $ws.Range("E2") "3"
Set-TextValue $ws.Range("F2") "24.284"
$ws.Range("G2").Value = "SEA"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "DAL"
$ws.Range("J2").Value = "W 24-13"
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = 0
